$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "51.753.32"
$ws.Range("E2").Value = "  +4.61%  "

# Row 3
$ws.Range("D3").Value = "2.771.16"
$ws.Range("E3").Value = "  +5.22%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "117.10"
$ws.Range("E5").Value = "  +4.55%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "333.21"

# Row 7
$ws.Range("E7").Value = "  +2.61%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.576"
$ws.Range("E9").Value = "  +6.05%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.25"
$ws.Range("E10").Value = "  +6.82%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0857"
$ws.Range("E11").Value = "  +5.84%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.27"
$ws.Range("E12").Value = "  +2.19%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.129"
$ws.Range("E13").Value = "  +2.03%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.66"
$ws.Range("E14").Value = "  +4.56%  "

# Row 15
$ws.Range("D15").Value = "3.207.97"

# Row 16
$ws.Range("D16").Value = "2.767.20"
$ws.Range("E16").Value = "  +5.17%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.890"
$ws.Range("E17").Value = "  +5.00%  "

# Row 18
$ws.Range("D18").Value = "51.670.79"
$ws.Range("E18").Value = "  +4.67%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.27"
$ws.Range("E19").Value = "  +11.28%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.56"
$ws.Range("E20").Value = "  +5.25%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.85"
$ws.Range("E21").Value = "  +2.45%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0974"
$ws.Range("E22").Value = "  +3.07%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "280.21"
$ws.Range("E23").Value = "  +3.78%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.02"
$ws.Range("E24").Value = "  +1.66%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.69"
$ws.Range("E25").Value = "  +6.48%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.87"
$ws.Range("E26").Value = "  +2.53%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.06%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.22"
$ws.Range("E28").Value = "  -0.37%  "

# Row 29
$ws.Range("E29").Value = "  +0.87%  "

# Row 30
$ws.Range("E30").Value = "  +3.61%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.27"
$ws.Range("E31").Value = "  +1.79%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.10"
$ws.Range("E32").Value = "  +1.28%  "

# Row 33
$ws.Range("E33").Value = "  +1.80%  "

# Row 34
$ws.Range("E34").Value = "  +1.02%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.18"
$ws.Range("E35").Value = "  +1.53%  "

# Row 36
$ws.Range("E36").Value = "  +0.01%  "

# Row 37
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.10"
$ws.Range("E37").Value = "  +2.92%  "

# Row 38
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.01"
$ws.Range("E38").Value = "  +2.62%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.27"
$ws.Range("E39").Value = "  +5.39%  "

# Row 40
$ws.Range("E40").Value = "  +10.26%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "127.70"
$ws.Range("E41").Value = "  -0.84%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "23.24"
$ws.Range("E42").Value = "  +5.69%  "

# Row 43
$ws.Range("E43").Value = "  +8.01%  "

# Row 44
$ws.Range("E44").Value = "  +2.93%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.47"
$ws.Range("E45").Value = "  +17.39%  "

# Row 46
$ws.Range("D46").Value = "2.087.23"
$ws.Range("E46").Value = "  +1.52%  "

# Row 47
$ws.Range("E47").Value = "  +3.71%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.25"
$ws.Range("E48").Value = "  +4.33%  "

# Row 49
$ws.Range("E49").Value = "  +6.51%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "60.81"
$ws.Range("E50").Value = "  +3.01%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.85"
$ws.Range("E51").Value = "  -0.59%  "
